$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the two new data rows (row 53 and 54) below the existing data.
# Reuse the existing date formatting (numFmtId 14, m/d/yyyy) from the cell
# above by copying its formats, rather than assigning a NumberFormat string
# (which would create a brand-new custom numFmt entry).
$ws.Range("A52").Copy()
$ws.Range("A53:A54").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A53").Value = 46029
$ws.Range("B53").Value = 1

$ws.Range("A54").Value = 46028
$ws.Range("B54").Value = 2

# Update the view to match the new selection/scroll position.
$ws.Range("A53:B54").Select()
$excel.ActiveWindow.ScrollRow = 40
